# Auto-generated script to apply scheduled-runner value updates
# to the Garuda_Profits Leve-profit tables across all 8 sheets.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1021.5
$ws.Range("I6").Value = 392.5
$ws.Range("J6").Value = 1650.5
$ws.Range("K6").Value = 1177.5
$ws.Range("L6").Value = 4951.5
$ws.Range("M6").Value = -1065.5
$ws.Range("N6").Value = -5175.5
$ws.Range("H28").Value = 127
$ws.Range("I28").Value = 127
$ws.Range("K28").Value = 127
$ws.Range("M28").Value = 358
$ws.Range("H33").Value = 408.8889
$ws.Range("I33").Value = 425.88235
$ws.Range("K33").Value = 425.88235
$ws.Range("M33").Value = -196.88235

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1067.1666
$ws.Range("I2").Value = 750
$ws.Range("J2").Value = 1384.3334
$ws.Range("K2").Value = 750
$ws.Range("L2").Value = 1384.3334
$ws.Range("M2").Value = -637
$ws.Range("N2").Value = -1610.3334
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H110").Value = 2223.5
$ws.Range("I110").Value = 2026.909
$ws.Range("J110").Value = 2944.3333
$ws.Range("K110").Value = 2026.909
$ws.Range("L110").Value = 2944.3333
$ws.Range("M110").Value = 18.09099999999989
$ws.Range("N110").Value = -7034.3333
$ws.Range("H116").Value = 1067.1666
$ws.Range("I116").Value = 750
$ws.Range("J116").Value = 1384.3334
$ws.Range("K116").Value = 750
$ws.Range("L116").Value = 1384.3334
$ws.Range("M116").Value = 1544
$ws.Range("N116").Value = -5972.3334

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1067.1666
$ws.Range("I3").Value = 750
$ws.Range("J3").Value = 1384.3334
$ws.Range("K3").Value = 750
$ws.Range("L3").Value = 1384.3334
$ws.Range("M3").Value = -636
$ws.Range("N3").Value = -1612.3334
$ws.Range("H63").Value = 49449.5
$ws.Range("J63").Value = 49449.5
$ws.Range("L63").Value = 49449.5
$ws.Range("N63").Value = -50821.5
$ws.Range("H66").Value = 49449.5
$ws.Range("J66").Value = 49449.5
$ws.Range("L66").Value = 148348.5
$ws.Range("N66").Value = -155212.5
$ws.Range("H94").Value = 1031.375
$ws.Range("I94").Value = 851.9286
$ws.Range("J94").Value = 1282.6
$ws.Range("K94").Value = 851.9286
$ws.Range("L94").Value = 1282.6
$ws.Range("M94").Value = -400.9286
$ws.Range("N94").Value = -2184.6
$ws.Range("H134").Value = 127189
$ws.Range("I134").Value = 251378
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 754134
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -751599
$ws.Range("N134").Value = -14070

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4653714
$ws.Range("I31").Value = 3439.4443
$ws.Range("J31").Value = 8001912
$ws.Range("K31").Value = 3439.4443
$ws.Range("L31").Value = 8001912
$ws.Range("M31").Value = -3144.4443
$ws.Range("N31").Value = -8002502
$ws.Range("H34").Value = 4653714
$ws.Range("I34").Value = 3439.4443
$ws.Range("J34").Value = 8001912
$ws.Range("K34").Value = 3439.4443
$ws.Range("L34").Value = 8001912
$ws.Range("M34").Value = -3237.4443
$ws.Range("N34").Value = -8002316
$ws.Range("H107").Value = 520.6667
$ws.Range("I107").Value = 513
$ws.Range("J107").Value = 545.2
$ws.Range("K107").Value = 513
$ws.Range("L107").Value = 545.2
$ws.Range("M107").Value = 1407
$ws.Range("N107").Value = -4385.2

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 503.97827
$ws.Range("I113").Value = 521.64703
$ws.Range("J113").Value = 493.6207
$ws.Range("K113").Value = 1564.94109
$ws.Range("L113").Value = 1480.8621
$ws.Range("M113").Value = 605.0589100000002
$ws.Range("N113").Value = -5820.8621
$ws.Range("H131").Value = 4489429.5
$ws.Range("J131").Value = 8008650.5
$ws.Range("L131").Value = 24025951.5
$ws.Range("N131").Value = -24036031.5

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 3547719.8
$ws.Range("I11").Value = 3767649.8
$ws.Range("J11").Value = 2668000
$ws.Range("K11").Value = 3767649.8
$ws.Range("L11").Value = 2668000
$ws.Range("M11").Value = -3767510.8
$ws.Range("N11").Value = -2668278
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H113").Value = 41667916
$ws.Range("I113").Value = 250000000
$ws.Range("K113").Value = 250000000
$ws.Range("M113").Value = -249997830

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value = 6612
$ws.Range("J76").Value = 13144
$ws.Range("L76").Value = 13144
$ws.Range("N76").Value = -13820
$ws.Range("H79").Value = 6612
$ws.Range("J79").Value = 13144
$ws.Range("L79").Value = 13144
$ws.Range("N79").Value = -15484
$ws.Range("H93").Value = 1734.25
$ws.Range("I93").Value = 1999.375
$ws.Range("J93").Value = 1204
$ws.Range("K93").Value = 1999.375
$ws.Range("L93").Value = 1204
$ws.Range("M93").Value = -751.375
$ws.Range("N93").Value = -3700
$ws.Range("H94").Value = 14000.5
$ws.Range("J94").Value = 14000.5
$ws.Range("L94").Value = 14000.5
$ws.Range("N94").Value = -15352.5
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 6002.5
$ws.Range("I14").Value = 2000
$ws.Range("J14").Value = 10005
$ws.Range("K14").Value = 2000
$ws.Range("L14").Value = 10005
$ws.Range("M14").Value = -1832
$ws.Range("N14").Value = -10341
$ws.Range("H69").Value = 19211.834
$ws.Range("J69").Value = 19211.834
$ws.Range("L69").Value = 19211.834
$ws.Range("N69").Value = -20709.834
$ws.Range("H72").Value = 19211.834
$ws.Range("J72").Value = 19211.834
$ws.Range("L72").Value = 57635.50199999999
$ws.Range("N72").Value = -65123.50199999999
$ws.Range("J76").Value = 19000
$ws.Range("L76").Value = 19000
$ws.Range("N76").Value = -19630
$ws.Range("J79").Value = 19000
$ws.Range("L79").Value = 19000
$ws.Range("N79").Value = -21184

